$wb = $excel.ActiveWorkbook

# --- Sheet "Activity Log - Part 1": scroll position changed (topLeftCell A58 -> A46) ---
$ws1 = $wb.Worksheets.Item("Activity Log - Part 1")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1

# --- Sheet "Activity Log - Part 2": new activity rows + updated times ---
$ws2 = $wb.Worksheets.Item("Activity Log - Part 2")
$ws2.Activate()

# Row 70 - existing entry, end time & note updated
$ws2.Range("E70").Value = 0.18055555555555555

# Order matters for shared-string table layout: G71's text must be interned
# before G70's so the underlying shared strings come out in the same order
# as the source workbook.
$ws2.Range("G71").Value = "Verified that LogicUnit behaves as expected for Timing Simultion. -DONE"
$ws2.Range("G70").Value = "Verified that timing simulation for ArithUnit and ExecUnit are still good. Communicated issue and fix to team members. Pushed changes to Github. -DONE"

# Row 71 - new entry
$ws2.Range("B71").Value = 6977
$ws2.Range("C71").Value = 43938
$ws2.Range("D71").Value = 0.18055555555555555
$ws2.Range("E71").Value = 0.1875

# Row 72 - new entry
$ws2.Range("B72").Value = 6977
$ws2.Range("C72").Value = 43938
$ws2.Range("D72").Value = 0.1875
$ws2.Range("E72").Value = 0.20833333333333334
$ws2.Range("G72").Value = "Discovered that test bench vector is ArithUnit01.tvs with team member. Changed it and replaced transcript files as waveforms are still valid. -DONE"

# Row 73 - new entry
$ws2.Range("B73").Value = 6977
$ws2.Range("C73").Value = 43938
$ws2.Range("D73").Value = 0.20833333333333334
$ws2.Range("E73").Value = 0.25
$ws2.Range("G73").Value = "Captured raw diagrams of timing waveforms from ModelSim for ExecUnit.vhd. Fixed a timing waveform in ArithUnit.vhd. -DONE"

# Active cell/selection moved to G74 on this sheet
$ws2.Range("G74").Select() | Out-Null
